# Apply the "complex scenario" test-data update to Sheet1.
#
# The original footer block (rows 24-26: footer line #1, a blank spacer
# row, footer line #2) needs to move down two rows (to rows 26-28) and a
# brand-new third footer row (29, column A only) gets appended. The two
# existing footer strings are reworded, and a third one is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old rows 24-26 down to rows 26-28 by inserting two fresh
# blank rows above the footer block (mirrors the row renumbering seen in
# the diff, including the pre-existing gap left by the missing row 23).
$ws.Rows.Item(24).Insert()
$ws.Rows.Item(24).Insert()

# Reword the (now shifted) footer lines and add the new third line.
$ws.Range("A26").Value = "Line#1 to drop"
$ws.Range("A28").Value = "Line#2 to drop"
$ws.Range("A29").Value = "Line#3 to drop"

# A29 is a brand-new row: give it the same look (fill/border) as the
# other "footer" cells in column A (e.g. A26) without touching B29:E29,
# which must stay absent, matching the target row 29 that only has a
# single populated cell.
$ws.Range("A26").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# Restore the selected cell to match the saved view state.
$ws.Range("A18").Select()
